$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 4648
$ws.Range("F3").Value = 2712
$ws.Range("F5").Value = 2719
$ws.Range("F9").Value = 1709
$ws.Range("F10").Value = 734
$ws.Range("F12").Value = 187
$ws.Range("F13").Value = 387
$ws.Range("F14").Value = 1054
$ws.Range("F15").Value = 292
$ws.Range("F17").Value = 64
$ws.Range("F18").Value = 516
$ws.Range("F22").Value = 734
$ws.Range("F23").Value = 140
$ws.Range("F24").Value = 28
$ws.Range("F25").Value = 493
$ws.Range("F26").Value = 1652
$ws.Range("F27").Value = 1427
$ws.Range("F28").Value = 303
$ws.Range("F29").Value = 38
$ws.Range("F30").Value = 1400
$ws.Range("F31").Value = 2260
$ws.Range("F32").Value = 369
$ws.Range("F33").Value = 22
$ws.Range("F34").Value = 594
$ws.Range("F35").Value = 111
$ws.Range("F36").Value = 51
$ws.Range("F38").Value = 762
$ws.Range("F39").Value = 1442
$ws.Range("F40").Value = 191
$ws.Range("F42").Value = 481
$ws.Range("F43").Value = 4
$ws.Range("F44").Value = 71

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F11").Value = 26

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 4648
$ws.Range("F3").Value = 2712
$ws.Range("F4").Value = 2719
$ws.Range("F5").Value = 1709
$ws.Range("F8").Value = 734
$ws.Range("F10").Value = 187
$ws.Range("F11").Value = 387
$ws.Range("F12").Value = 1054
$ws.Range("F13").Value = 292
$ws.Range("F15").Value = 64
$ws.Range("F16").Value = 516
$ws.Range("F19").Value = 734
$ws.Range("F20").Value = 140
$ws.Range("F24").Value = 28
$ws.Range("F25").Value = 493
$ws.Range("F26").Value = 1652
$ws.Range("F27").Value = 1427
$ws.Range("F28").Value = 303
$ws.Range("F29").Value = 38
$ws.Range("F32").Value = 2260
$ws.Range("F33").Value = 369
$ws.Range("F36").Value = 22
$ws.Range("F37").Value = 26
$ws.Range("F38").Value = 594
$ws.Range("F39").Value = 111
$ws.Range("F40").Value = 51
$ws.Range("F42").Value = 762
$ws.Range("F43").Value = 1442
$ws.Range("F45").Value = 191
$ws.Range("F46").Value = 481
$ws.Range("F47").Value = 71
